# "Add files via upload" — append the latest expense-record entries to the
# "Rafid" expense tracker (Money Manager Python / expense_record.xlsx) and
# bring the existing "Date" column onto the newer yyyy-mm-dd hh:mm:ss
# display format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 19: a same-day Food expense (date-only serial) --------------
$ws.Range("A19").Value = 45419
$ws.Range("B19").Value = "Food"
$ws.Range("C19").Value = 1000

# --- New row 20: a Transport expense logged later the same day, with a
#     time-of-day component in the serial value ---------------------------
$ws.Range("A20").Value = 45419.96818364583
$ws.Range("B20").Value = "Transport"
$ws.Range("C20").Value = 100

# --- New row 21: a Food expense whose date was typed/imported as plain
#     text ("2024-05-07") rather than a real date serial -------------------
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "2024-05-07"
$ws.Range("B21").Value = "Food"
$ws.Range("C21").Value = 10

# --- Refresh the Date column's display format for every dated row
#     (existing history A2:A18 plus the two new date rows A19:A20) so they
#     all render as "YYYY-MM-DD HH:MM:SS" ----------------------------------
$ws.Range("A2:A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
